$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3847390.2
$ws.Range("J17").Value = 3847390.2
$ws.Range("L17").Value = 11542170.6
$ws.Range("N17").Value = -11542506.6

$ws.Range("H40").Value = 4134.2085
$ws.Range("I40").Value = 2764
$ws.Range("J40").Value = 6874.625
$ws.Range("K40").Value = 2764
$ws.Range("L40").Value = 6874.625
$ws.Range("M40").Value = -2589
$ws.Range("N40").Value = -7224.625

$ws.Range("H64").Value = 7309.852
$ws.Range("I64").Value = 3397.9375
$ws.Range("K64").Value = 3397.9375
$ws.Range("M64").Value = -3149.9375

$ws.Range("H67").Value = 7309.852
$ws.Range("I67").Value = 3397.9375
$ws.Range("K67").Value = 3397.9375
$ws.Range("M67").Value = -2539.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 13995
$ws.Range("I61").Value = 14667.917
$ws.Range("K61").Value = 14667.917
$ws.Range("M61").Value = -14455.917

$ws.Range("H94").Value = 50250
$ws.Range("J94").Value = 50250
$ws.Range("L94").Value = 50250
$ws.Range("N94").Value = -52052

$ws.Range("H136").Value = 13995
$ws.Range("I136").Value = 14667.917
$ws.Range("K136").Value = 44003.751
$ws.Range("M136").Value = -41453.751

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 4998
$ws.Range("I7").Value = 4996
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 4996
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -4883
$ws.Range("N7").Value = -5226

$ws.Range("H14").Value = 10000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 10000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -10344

$ws.Range("I20").Value = 4282.2
$ws.Range("J20").Value = 5587.4546
$ws.Range("K20").Value = 4282.2
$ws.Range("L20").Value = 5587.4546
$ws.Range("M20").Value = -4035.2
$ws.Range("N20").Value = -6081.4546

$ws.Range("H22").Value = 371.1875
$ws.Range("I22").Value = 211.58333
$ws.Range("J22").Value = 850
$ws.Range("K22").Value = 211.58333
$ws.Range("L22").Value = 850
$ws.Range("M22").Value = -38.58332999999999
$ws.Range("N22").Value = -1196

$ws.Range("H86").Value = 2394
$ws.Range("I86").Value = 2106.625
$ws.Range("J86").Value = 3773.4
$ws.Range("K86").Value = 2106.625
$ws.Range("L86").Value = 3773.4
$ws.Range("M86").Value = -983.625
$ws.Range("N86").Value = -6019.4

$ws.Range("H89").Value = 2394
$ws.Range("I89").Value = 2106.625
$ws.Range("J89").Value = 3773.4
$ws.Range("K89").Value = 10533.125
$ws.Range("L89").Value = 18867
$ws.Range("M89").Value = -4917.125
$ws.Range("N89").Value = -30099

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 26100
$ws.Range("J50").Value = 26100
$ws.Range("L50").Value = 26100
$ws.Range("N50").Value = -27350

$ws.Range("H80").Value = 39405.89
$ws.Range("J80").Value = 39405.89
$ws.Range("L80").Value = 39405.89
$ws.Range("N80").Value = -41651.89

$ws.Range("H83").Value = 39405.89
$ws.Range("J83").Value = 39405.89
$ws.Range("L83").Value = 118217.67
$ws.Range("N83").Value = -129449.67

$ws.Range("H86").Value = 48993.418
$ws.Range("I86").Value = 6854.2
$ws.Range("K86").Value = 6854.2
$ws.Range("M86").Value = -5731.2

$ws.Range("H89").Value = 48993.418
$ws.Range("I89").Value = 6854.2
$ws.Range("K89").Value = 34271
$ws.Range("M89").Value = -28655

$ws.Range("H99").Value = 3216.4167
$ws.Range("J99").Value = 2273.8572
$ws.Range("L99").Value = 2273.8572
$ws.Range("N99").Value = -5269.8572

$ws.Range("H126").Value = 3216.4167
$ws.Range("J126").Value = 2273.8572
$ws.Range("L126").Value = 6821.571599999999
$ws.Range("N126").Value = -11761.5716

$ws.Range("H132").Value = 5540.125
$ws.Range("I132").Value = 5578.2
$ws.Range("K132").Value = 16734.6
$ws.Range("M132").Value = -14204.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 947.2222
$ws.Range("J12").Value = 1217.5714
$ws.Range("L12").Value = 3652.7142
$ws.Range("N12").Value = -3998.7142

$ws.Range("H92").Value = 249.4
$ws.Range("I92").Value = 350
$ws.Range("K92").Value = 1050
$ws.Range("M92").Value = 198

$ws.Range("H98").Value = 303
$ws.Range("I98").Value = 303
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 909
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 589
$ws.Range("N98").ClearContents()

$ws.Range("H107").Value = 398.30768
$ws.Range("I107").Value = 118
$ws.Range("J107").Value = 449.27274
$ws.Range("K107").Value = 354
$ws.Range("L107").Value = 1347.81822
$ws.Range("M107").Value = 1566
$ws.Range("N107").Value = -5187.81822

$ws.Range("H113").Value = 1007.6
$ws.Range("J113").Value = 1072.5
$ws.Range("L113").Value = 3217.5
$ws.Range("N113").Value = -7557.5

$ws.Range("H132").Value = 2023.421
$ws.Range("I132").Value = 1806.125
$ws.Range("J132").Value = 2181.4546
$ws.Range("K132").Value = 16255.125
$ws.Range("L132").Value = 19633.0914
$ws.Range("M132").Value = -13725.125
$ws.Range("N132").Value = -24693.0914

$ws.Range("H140").Value = 2861.25
$ws.Range("I140").Value = 2649.1667
$ws.Range("K140").Value = 7947.500100000001
$ws.Range("M140").Value = -2767.500100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15260.223
$ws.Range("I70").Value = 13593.444
$ws.Range("J70").Value = 16927
$ws.Range("K70").Value = 13593.444
$ws.Range("L70").Value = 16927
$ws.Range("M70").Value = -13323.444
$ws.Range("N70").Value = -17467

$ws.Range("H73").Value = 15260.223
$ws.Range("I73").Value = 13593.444
$ws.Range("J73").Value = 16927
$ws.Range("K73").Value = 13593.444
$ws.Range("L73").Value = 16927
$ws.Range("M73").Value = -12657.444
$ws.Range("N73").Value = -18799

$ws.Range("H93").Value = 37027.5
$ws.Range("J93").Value = 37027.5
$ws.Range("L93").Value = 37027.5
$ws.Range("N93").Value = -40771.5

$ws.Range("H132").Value = 1626.6666
$ws.Range("I132").Value = 1190
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 3570
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -1040
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1808.3529
$ws.Range("I16").Value = 1667.2858
$ws.Range("K16").Value = 1667.2858
$ws.Range("M16").Value = -1497.2858

$ws.Range("H22").Value = 1089.6923
$ws.Range("J22").Value = 1076
$ws.Range("L22").Value = 1076
$ws.Range("N22").Value = -1666

$ws.Range("H27").Value = 1089.6923
$ws.Range("J27").Value = 1076
$ws.Range("L27").Value = 1076
$ws.Range("N27").Value = -1290

$ws.Range("H40").Value = 2202.3333
$ws.Range("I40").Value = 2052.5
$ws.Range("J40").Value = 5199
$ws.Range("K40").Value = 2052.5
$ws.Range("L40").Value = 5199
$ws.Range("M40").Value = -1916.5
$ws.Range("N40").Value = -5471

$ws.Range("H55").Value = 412.16666
$ws.Range("I55").Value = 427.3889
$ws.Range("K55").Value = 427.3889
$ws.Range("M55").Value = -254.3889

$ws.Range("H122").Value = 4586.0835
$ws.Range("J122").Value = 4999.25
$ws.Range("L122").Value = 14997.75
$ws.Range("N122").Value = -19897.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 10000
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H100").Value = 612.1429000000001
$ws.Range("I100").Value = 810.25
$ws.Range("J100").Value = 348
$ws.Range("K100").Value = 1620.5
$ws.Range("L100").Value = 696
$ws.Range("M100").Value = -1079.5
$ws.Range("N100").Value = -1778

$ws.Range("H136").Value = 4629.9
$ws.Range("I136").Value = 2672.1667
$ws.Range("K136").Value = 8016.500100000001
$ws.Range("M136").Value = -5466.500100000001
